$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.015.76"
$ws.Range("E2").Value = "  -4.13%  "
$ws.Range("D3").Value = "3.418.07"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'237.41"
$ws.Range("D6").Value = "'640.00"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  -3.44%  "
$ws.Range("D8").Value = "'0.401"
$ws.Range("E8").Value = "  -5.86%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "'0.966"
$ws.Range("E10").Value = "  -6.57%  "
$ws.Range("D11").Value = "3.417.04"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  -4.63%  "
$ws.Range("D13").Value = "'41.47"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").Value = "'6.21"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "93.817.41"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("D16").Value = "4.058.80"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'0.0000250"
$ws.Range("E17").Value = "  -3.01%  "
$ws.Range("D18").Value = "'8.29"
$ws.Range("E18").Value = "  -10.21%  "
$ws.Range("D19").Value = "3.417.43"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "'17.41"
$ws.Range("E20").Value = "  -3.46%  "
$ws.Range("D21").Value = "'11.57"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'498.34"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("B23").Value = "Stellar"
$ws.Range("C23").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D23").Value = "'0.469"
$ws.Range("E23").Value = "  -10.19%  "
$ws.Range("D24").Value = "'3.23"
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("D25").Value = "'0.0000192"
$ws.Range("E25").Value = "  -4.64%  "
$ws.Range("D26").Value = "'6.49"
$ws.Range("E26").Value = "  -5.66%  "
$ws.Range("D27").Value = "'91.02"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").Value = "3.601.41"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "'11.88"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").Value = "'11.49"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("D33").Value = "'0.136"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'0.175"
$ws.Range("E35").Value = "  -8.12%  "
$ws.Range("D36").Value = "'29.47"
$ws.Range("E36").Value = "  +2.39%  "
$ws.Range("D37").Value = "'0.547"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").Value = "'541.66"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("D39").Value = "'7.67"
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("E40").Value = "  -2.87%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'0.150"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "'0.904"
$ws.Range("E43").Value = "  +6.57%  "
$ws.Range("D44").Value = "'24.00"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "'3.72"
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").Value = "'1.71"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").Value = "'5.63"
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "'2.21"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'3.29"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").Value = "'0.0405"
$ws.Range("E50").Value = "  -4.98%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'54.33"
$ws.Range("E51").Value = "  -2.04%  "
